# Registree stats backup on Wed 21 Apr 2021 09:29:43 SAST
# Refresh the four timestamped report headers, a handful of per-registree
# "Voter" flags that flipped between the 09:27 and 09:29 snapshots, the
# derived "Number of voters" totals, and the 410W voting-by-club table
# (Durbanville dropped out, Cape Of Good Hope gained a voter).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MD410 Attendance
# ---------------------------------------------------------------------
$wsMD410 = $wb.Worksheets.Item("MD410 Attendance")
$wsMD410.Range("A1").Value = "MD410 Registrees as of 21/04/2021 09:29"
$wsMD410.Range("A239").Value = "Number of voters: 99"

# ---------------------------------------------------------------------
# 410E Attendance
# ---------------------------------------------------------------------
$wsE = $wb.Worksheets.Item("410E Attendance")
$wsE.Range("A1").Value = "410E Registrees as of 21/04/2021 09:29"
# Tracey Polkinghorne (row 81) is no longer recorded as a voter.
$wsE.Range("E81").Value = "No"

# ---------------------------------------------------------------------
# 410W Attendance
# ---------------------------------------------------------------------
$wsW = $wb.Worksheets.Item("410W Attendance")
$wsW.Range("A1").Value = "410W Registrees as of 21/04/2021 09:29"
# Nelle Buhrmann (row 13) is now a voter.
$wsW.Range("E13").Value = "Yes"
# Dawid Jacobs (row 55) is no longer a voter.
$wsW.Range("E55").Value = "No"
# Yvonne Maureen Poyowe (row 101) is now a voter.
$wsW.Range("E101").Value = "Yes"
$wsW.Range("A152").Value = "Number of voters: 56"

# ---------------------------------------------------------------------
# 410E Voting
# ---------------------------------------------------------------------
$wsEV = $wb.Worksheets.Item("410E Voting")
$wsEV.Range("A1").Value = "410E Voting details as of 21/04/2021 09:29"

# ---------------------------------------------------------------------
# 410W Voting
# ---------------------------------------------------------------------
$wsWV = $wb.Worksheets.Item("410W Voting")
$wsWV.Range("A1").Value = "410W Voting details as of 21/04/2021 09:29"
# Cape Of Good Hope (row 5) gained a voter.
$wsWV.Range("B5").Value = 3
# Durbanville (row 9) dropped out of the voting-club list entirely; deleting
# its row shifts every following club up by one and lets Excel keep the
# dimension/footer rows in sync automatically.
$wsWV.Rows.Item(9).Delete()
# Footer totals: one fewer club, one more voter overall.
$wsWV.Range("A36").Value = "Number of clubs: 33"
$wsWV.Range("A37").Value = "Number of voters: 56"
